$d = $word.ActiveDocument

# The document ends with a paragraph that only contains a single space
# run. The edit inserts five new paragraphs of notes text, two blank
# paragraphs, and then re-applies "en-GB" language formatting (on both
# the paragraph mark and the run) to that trailing paragraph.

$target = $d.Paragraphs.Last

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pPrLang = "<w:pPr><w:rPr><w:lang w:val=`"en-GB`"/></w:rPr></w:pPr>"
$rPrLang = "<w:rPr><w:lang w:val=`"en-GB`"/></w:rPr>"

function NewPara($text) {
    if ($text -eq $null) {
        return "<w:p $ns>$pPrLang</w:p>"
    }
    return "<w:p $ns>$pPrLang<w:r>$rPrLang<w:t xml:space=`"preserve`">$text</w:t></w:r></w:p>"
}

$xml = ""
$xml += NewPara("Inputs size: Number of feature dimensions (e.g., number of data channels). Often set to 1.")
$xml += NewPara("Sequence length: Length of data snippet used for learning.")
$xml += NewPara("Batch size: Number of sequences to input per epoch.")
$xml += NewPara("Hidden size: Number of units in the hidden layer.")
$xml += NewPara("Number of layers: Number of hidden layers stacked on each other.")
$xml += NewPara($null)
$xml += NewPara($null)
# Re-create the trailing paragraph itself (same single-space run) so it
# picks up the en-GB paragraph-mark + run language formatting.
$xml += NewPara(" ")

$full = $d.Range($target.Range.Start, $target.Range.End)
$full.InsertXML($xml)
